$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Header contact line: split the single centered line into two centered
#    lines - "Pittsburgh, PA ... gmail.com" stays, and the web links move to
#    their own paragraph right below it.
# ---------------------------------------------------------------------------
$pContact = $d.Paragraphs.Item(2)
$pContact.Range.Find.Execute(
    " | slimeq.github.io | github.com/SlimeQ | linkedin.com/in/quincy-campbell-131559b2",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

$pContact.Range.InsertParagraphAfter()
$pLinks = $d.Paragraphs.Item($pContact.Index + 1)
$pLinks.Range.Text = "slimeq.github.io | github.com/SlimeQ | linkedin.com/in/quincy-campbell-131559b2"

# ---------------------------------------------------------------------------
# 2) Section-header "space before" tightened from 10pt (200 twips) to
#    6pt (120 twips) throughout the document: SUMMARY, PRIMARY SKILLS,
#    PROFESSIONAL EXPERIENCE, EDUCATION, SELECTED WORK / PORTFOLIO.
# ---------------------------------------------------------------------------
$headerTitles = @("SUMMARY", "PRIMARY SKILLS", "PROFESSIONAL EXPERIENCE", "EDUCATION", "SELECTED WORK / PORTFOLIO")
foreach ($title in $headerTitles) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $pp = $d.Paragraphs.Item($i)
        $txt = $pp.Range.Text.TrimEnd([char]13)
        if ($txt -eq $title) {
            $pp.SpaceBefore = 6
            break
        }
    }
}

# ---------------------------------------------------------------------------
# 3) PRIMARY SKILLS: turn the single bullet-separated paragraph into six
#    separate "List Bullet" paragraphs.
# ---------------------------------------------------------------------------
$pSkills = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    $txt = $pp.Range.Text.TrimEnd([char]13)
    if ($txt.StartsWith("Unity / C# (architecture")) {
        $pSkills = $pp
        break
    }
}

# Trim the paragraph down to just the first bullet's text (this keeps the
# run's existing rPr intact instead of rebuilding it from scratch).
$pSkills.Range.Find.Execute(
    " " + [char]0x2022 + " Deployments: Windows, WebGL, Android, iOS " + [char]0x2022 + " XR: Magic Leap 2 " + [char]0x2022 + " Integration: REST APIs, WebSockets, networking (Mirror / Netcode for GameObjects) " + [char]0x2022 + " Graphics: shaders/materials, procedural mesh, 3D math/physics " + [char]0x2022 + " Backend: ASP.NET Core, PostgreSQL",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# Apply the List Bullet style + tight after-spacing (matches the rest of the
# bulleted sections in the resume).
$pSkills.Format.Style = "List Bullet"
$pSkills.SpaceAfter = 1

# Re-apply direct character formatting to the run text only (exclude the
# trailing paragraph mark so Word doesn't also stamp mark-run-properties
# onto pPr).
$rFix = $d.Range($pSkills.Range.Start, $pSkills.Range.End - 1)
$rFix.Font.Name = "Calibri"
$rFix.Font.Size = 10

$restOfSkills = @(
    "Deployments: Windows, WebGL, Android, iOS",
    "XR: Magic Leap 2 (Android), VR/AR experience",
    "Integration: REST APIs, WebSockets, analytics/data export patterns",
    "Graphics: shaders/materials, procedural mesh, 3D math/physics",
    "Backend: ASP.NET Core / Blazor, PostgreSQL"
)

$prev = $pSkills
foreach ($line in $restOfSkills) {
    $prev.Range.InsertParagraphAfter()
    $np = $d.Paragraphs.Item($prev.Index + 1)
    $np.Range.Text = $line
    $prev = $np
}

Write-Output "done"
